$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 4.739369273631283
$ws.Range("C2").Value = 4.733806412751601
$ws.Range("D2").Value = 459.0342543959424
$ws.Range("E2").Value = 0.01137594362591159
$ws.Range("B3").Value = 4.739369273631283
$ws.Range("C3").Value = 4.733806412751601
$ws.Range("D3").Value = 459.0342543959424
$ws.Range("E3").Value = 0.01072878202702245
$ws.Range("B4").Value = 4.977612155374525
$ws.Range("C4").Value = 4.958776757284687
$ws.Range("D4").Value = 7.247728239222934
$ws.Range("E4").Value = 0.008927012830482495
$ws.Range("B5").Value = 4.977612155374525
$ws.Range("C5").Value = 4.958776757284687
$ws.Range("D5").Value = 7.247728239222934
$ws.Range("E5").Value = 0.005790581315374009
$ws.Range("B6").Value = 4.977612155374525
$ws.Range("C6").Value = 4.958776757284687
$ws.Range("D6").Value = 7.247728239222934
$ws.Range("E6").Value = 0.003464691520643704
$ws.Range("B7").Value = 4.977612155374525
$ws.Range("C7").Value = 4.958776757284687
$ws.Range("D7").Value = 7.247728239222934
$ws.Range("E7").Value = 0.009663225338040772
$ws.Range("B8").Value = 4.968584659228433
$ws.Range("C8").Value = 4.988778195936124
$ws.Range("D8").Value = 3.661022268609827
$ws.Range("E8").Value = 0.00990851032948474
$ws.Range("B9").Value = 4.975250419765318
$ws.Range("C9").Value = 4.993348893247816
$ws.Range("D9").Value = 2.15872358669173
$ws.Range("E9").Value = 0.009580593045533901
$ws.Range("B10").Value = 4.975250419765318
$ws.Range("C10").Value = 4.993348893247816
$ws.Range("D10").Value = 2.15872358669173
$ws.Range("E10").Value = 0.01096138785588309
$ws.Range("B11").Value = 4.975250419765318
$ws.Range("C11").Value = 4.993348893247816
$ws.Range("D11").Value = 2.15872358669173
$ws.Range("E11").Value = 0.01052400080951986
$ws.Range("B12").Value = 4.986390794584488
$ws.Range("C12").Value = 5.014277080453552
$ws.Range("D12").Value = 1.293141969754038
$ws.Range("E12").Value = 0.008854183951754073
$ws.Range("B13").Value = 4.993383311261053
$ws.Range("C13").Value = 4.993348893247816
$ws.Range("D13").Value = 0.2879137868798451
$ws.Range("E13").Value = 0.01151609599376997
$ws.Range("B14").Value = 4.991681653200199
$ws.Range("C14").Value = 4.998228616194902
$ws.Range("D14").Value = 0.2380379360960141
$ws.Range("E14").Value = 0.01092924081609513
$ws.Range("B15").Value = 5.000250444920759
$ws.Range("C15").Value = 5.000280965979282
$ws.Range("D15").Value = 0.005935663731419322
$ws.Range("E15").Value = 0.006993530024554418
$ws.Range("B16").Value = 4.999967104180334
$ws.Range("C16").Value = 5.00025966044814
$ws.Range("D16").Value = 0.005497475709185652
$ws.Range("E16").Value = 0.007351238293140759
$ws.Range("B17").Value = 4.999967104180334
$ws.Range("C17").Value = 5.00025966044814
$ws.Range("D17").Value = 0.005497475709185652
$ws.Range("E17").Value = 0.008268206027352856
$ws.Range("B18").Value = 4.999967104180334
$ws.Range("C18").Value = 5.00025966044814
$ws.Range("D18").Value = 0.005497475709185652
$ws.Range("E18").Value = 0.004333873690907626
$ws.Range("B19").Value = 4.999967104180334
$ws.Range("C19").Value = 5.00025966044814
$ws.Range("D19").Value = 0.005497475709185652
$ws.Range("E19").Value = 0.006154916302306768
$ws.Range("B20").Value = 5.000126848365561
$ws.Range("C20").Value = 4.99974426857862
$ws.Range("D20").Value = 0.005311591277877292
$ws.Range("E20").Value = 0.004474406245892526
$ws.Range("B21").Value = 4.999978002107094
$ws.Range("C21").Value = 4.999697872214196
$ws.Range("D21").Value = 0.005219978598075844
$ws.Range("E21").Value = 0.01207370020522915
$ws.Range("B22").Value = 4.999970928302389
$ws.Range("C22").Value = 5.000062526917813
$ws.Range("D22").Value = 0.005162089461696572
$ws.Range("E22").Value = 0.01612764977485064
$ws.Range("B23").Value = 4.999930523798098
$ws.Range("C23").Value = 4.999957042482018
$ws.Range("D23").Value = 0.005074561432837549
$ws.Range("E23").Value = 0.03004136544881693
$ws.Range("B24").Value = 4.999868308480597
$ws.Range("C24").Value = 4.999892108024093
$ws.Range("D24").Value = 0.005066590317342493
$ws.Range("E24").Value = 0.2148146314007081
$ws.Range("B25").Value = 4.99993827929849
$ws.Range("C25").Value = 4.999905950516123
$ws.Range("D25").Value = 0.005066476609897517
$ws.Range("E25").Value = 0.5542038037109747
$ws.Range("B26").Value = 4.999893700006857
$ws.Range("C26").Value = 4.999897131289655
$ws.Range("D26").Value = 0.00506265596631361
$ws.Range("E26").Value = 2.014531456672368
